$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K2").Value = "hits[0].id||hits[1].id"
$ws.Range("L2:L3").Select()
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 1
